$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-of date range) ---
$ws.Range("A8").Value = "Volume 30   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/18/2023  Through  9/24/2023"

# --- Straightforward numeric value updates (style/type unchanged) ---
$ws.Range("N14").Value = -86.95652173913
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -60
$ws.Range("I15").Value = 20
$ws.Range("J15").Value = 32
$ws.Range("K15").Value = -37.5
$ws.Range("L15").Value = 11.111111111111
$ws.Range("M15").Value = 5.263157894736
$ws.Range("N15").Value = -66.666666666666
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 11.111111111111
$ws.Range("I16").Value = 104
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 4
$ws.Range("L16").Value = -2.803738317757
$ws.Range("M16").Value = -58.565737051792
$ws.Range("N16").Value = -87.317073170731
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 9.090909090909
$ws.Range("I17").Value = 256
$ws.Range("J17").Value = 298
$ws.Range("K17").Value = -14.093959731543
$ws.Range("L17").Value = -14.093959731543
$ws.Range("M17").Value = -0.389105058365
$ws.Range("N17").Value = -52.416356877323
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -85.714285714285
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -54.545454545454
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 95
$ws.Range("K18").Value = -41.052631578947
$ws.Range("L18").Value = -37.777777777777
$ws.Range("M18").Value = -79.56204379562
$ws.Range("N18").Value = -92.911392405063
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = -30.30303030303
$ws.Range("I19").Value = 299
$ws.Range("J19").Value = 325
$ws.Range("K19").Value = -8
$ws.Range("L19").Value = 20.564516129032
$ws.Range("M19").Value = -32.808988764044
$ws.Range("N19").Value = -90.401284109149
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -18.181818181818
$ws.Range("I20").Value = 154
$ws.Range("J20").Value = 181
$ws.Range("K20").Value = -14.917127071823
$ws.Range("L20").Value = 19.37984496124
$ws.Range("M20").Value = -22.613065326633
$ws.Range("N20").Value = -87.35632183908
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -32.258064516129
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 102
$ws.Range("H21").Value = -19.607843137254
$ws.Range("I21").Value = 892
$ws.Range("J21").Value = 1039
$ws.Range("K21").Value = -14.14821944177
$ws.Range("L21").Value = -0.888888888888
$ws.Range("M21").Value = -38.904109589041
$ws.Range("N21").Value = -86.41072516758
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -40.74074074074
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = 7.865168539325
$ws.Range("I24").Value = 868
$ws.Range("J24").Value = 956
$ws.Range("K24").Value = -9.205020920502
$ws.Range("L24").Value = 29.166666666666
$ws.Range("M24").Value = 14.060446780551
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = 38.461538461538
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 47
$ws.Range("H25").Value = 23.404255319148
$ws.Range("I25").Value = 496
$ws.Range("J25").Value = 391
$ws.Range("K25").Value = 26.854219948849
$ws.Range("L25").Value = 42.120343839541
$ws.Range("M25").Value = -17.744610281923
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 33
$ws.Range("J26").Value = 43
$ws.Range("K26").Value = -23.255813953488
$ws.Range("L26").Value = -5.714285714285
$ws.Range("J27").Value = 39
$ws.Range("K27").Value = -5.128205128205
$ws.Range("L27").Value = 0
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -100
$ws.Range("L28").Value = -64.444444444444
$ws.Range("N28").Value = -85.046728971962
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -100
$ws.Range("L29").Value = -68.571428571428
$ws.Range("N29").Value = -88.421052631578

# --- Cells that flip from a numeric count to the text placeholder "0" ---
# (used by the report whenever a 28-day count is zero)
$donorText = $ws.Range("C22")
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$donorText.Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$donorText.Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$donorText.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$donorText.Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$donorText.Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$donorText.Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$donorText.Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$donorText.Copy()
$ws.Range("F29").PasteSpecial(-4122)

# --- Cells that flip from the text placeholder to a real numeric value ---
$donorNum15 = $ws.Range("F18")
$ws.Range("C15").Value = 1
$donorNum15.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D27").Value = 2
$donorNum15.Copy()
$ws.Range("D27").PasteSpecial(-4122)

$donorNum16 = $ws.Range("L26")
$ws.Range("E27").Value = -100
$donorNum16.Copy()
$ws.Range("E27").PasteSpecial(-4122)

$excel.CutCopyMode = $false
